{"js": "// Locate the 5-paragraph block that starts with \"Resource.apply(Resource res) : Resource.\"\n// and rewrite its text content in place (paragraph count / formatting untouched):\n//\n//   1. \"Resource.apply(Resource res) : Resource.\"\n//   2. \"Transform. Update. Apply to player context resource holds for all occurrences.\"\n//   3. \"\"  (empty)\n//   4. \"Resource.query(Resource pattern);\"\n//   5. \"Apply to Model quad pattern performs resource activation (transform result).\"\n//\n// becomes\n//\n//   1. \"Resource.apply(Resource pattern) : Resource. Transform. Update. Apply pattern\n//       query / match: add / modify corresponding occurrences to player context resource.\"\n//   2. \"\"  (empty)\n//   3. \"Resource.query(Resource pattern); Apply to Model. Quad pattern matches. If none\n//       then build Resource from monadic resource factory. Performs resource activation\n//       (messages transform results, apply occurrences).\"\n//   4. \"\"  (empty)\n//   5. \"Metamodel messages: (match, apply) CSPO quads for each Resource hierarchy new\n//       instance: quads message. Apply occurrences to each local matching CSPO.\n//       Metamodels aggregate new occurrences.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Resource.apply(Resource res) : Resource.\";\nlet startIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === marker) {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1) {\n  throw new Error(\"Could not find the target paragraph block.\");\n}\n\nconst newTexts = [\n  \"Resource.apply(Resource pattern) : Resource. Transform. Update. Apply pattern query / match: add / modify corresponding occurrences to player context resource.\",\n  \"\",\n  \"Resource.query(Resource pattern); Apply to Model. Quad pattern matches. If none then build Resource from monadic resource factory. Performs resource activation (messages transform results, apply occurrences).\",\n  \"\",\n  \"Metamodel messages: (match, apply) CSPO quads for each Resource hierarchy new instance: quads message. Apply occurrences to each local matching CSPO. Metamodels aggregate new occurrences.\"\n];\n\nfor (let i = 0; i < newTexts.length; i++) {\n  const paragraph = paragraphs.items[startIndex + i];\n  // Replacing the paragraph's whole range with the new string (even an empty\n  // one) swaps the run's text while leaving the paragraph itself - and its\n  // run/paragraph formatting - in place.\n  paragraph.insertText(newTexts[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Locate the 5-paragraph block that starts with \"Resource.apply(Resource res) : Resource.\"\n# and rewrite its text content in place (paragraph count / formatting untouched):\n#\n#   1. \"Resource.apply(Resource res) : Resource.\"\n#   2. \"Transform. Update. Apply to player context resource holds for all occurrences.\"\n#   3. \"\"  (empty)\n#   4. \"Resource.query(Resource pattern);\"\n#   5. \"Apply to Model quad pattern performs resource activation (transform result).\"\n#\n# becomes\n#\n#   1. \"Resource.apply(Resource pattern) : Resource. Transform. Update. Apply pattern\n#       query / match: add / modify corresponding occurrences to player context resource.\"\n#   2. \"\"  (empty)\n#   3. \"Resource.query(Resource pattern); Apply to Model. Quad pattern matches. If none\n#       then build Resource from monadic resource factory. Performs resource activation\n#       (messages transform results, apply occurrences).\"\n#   4. \"\"  (empty)\n#   5. \"Metamodel messages: (match, apply) CSPO quads for each Resource hierarchy new\n#       instance: quads message. Apply occurrences to each local matching CSPO.\n#       Metamodels aggregate new occurrences.\"\n\n$d = $word.ActiveDocument\n\n$marker = \"Resource.apply(Resource res) : Resource.\"\n$startIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\")\n    if ($text -eq $marker) {\n        $startIndex = $i\n        break\n    }\n}\n\nif ($startIndex -eq -1) {\n    throw \"Could not find the target paragraph block.\"\n}\n\n$newTexts = @(\n    \"Resource.apply(Resource pattern) : Resource. Transform. Update. Apply pattern query / match: add / modify corresponding occurrences to player context resource.\",\n    \"\",\n    \"Resource.query(Resource pattern); Apply to Model. Quad pattern matches. If none then build Resource from monadic resource factory. Performs resource activation (messages transform results, apply occurrences).\",\n    \"\",\n    \"Metamodel messages: (match, apply) CSPO quads for each Resource hierarchy new instance: quads message. Apply occurrences to each local matching CSPO. Metamodels aggregate new occurrences.\"\n)\n\nfor ($j = 0; $j -lt $newTexts.Length; $j++) {\n    $p = $d.Paragraphs.Item($startIndex + $j)\n    # Assigning .Range.Text swaps the run's text while leaving the paragraph\n    # mark - and the paragraph/run formatting - untouched.\n    $p.Range.Text = $newTexts[$j]\n}\n"}
